# Update the 'Datos actualizados' timestamp in cell A1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = 'Datos actualizados a 28 de Abril de 2020 a las 08:52'

# Refresh the country ranking table (rows re-sorted by "Casos totales" after new data came in)
$updates = @(
    @{ Row = 39; Country = 'Ucrania'; Values = @(9410, 401, 992, 8179, 121, 19, 239) }
    @{ Row = 40; Country = 'Indonesia'; Values = @(9096, 0, 1151, 7180, 0, 0, 765) }
    @{ Row = 45; Country = 'Chequia'; Values = @(7449, 4, 2842, 4384, 73, 0, 223) }
    @{ Row = 99; Country = 'Kirguistan'; Values = @(708, 13, 416, 284, 13, 0, 8) }
    @{ Row = 100; Country = 'Honduras'; Values = @(702, 41, 79, 559, 10, 3, 64) }
    @{ Row = 101; Country = 'Niger'; Values = @(701, 0, 385, 287, 0, 0, 29) }
    @{ Row = 102; Country = 'Costa Rica'; Values = @(697, 0, 287, 404, 8, 0, 6) }
    @{ Row = 119; Country = 'El Salvador'; Values = @(345, 22, 97, 240, 4, 0, 8) }
    @{ Row = 120; Country = 'Estado de Palestina'; Values = @(342, 0, 83, 257, 0, 0, 2) }
    @{ Row = 121; Country = 'Mauricio'; Values = @(334, 0, 302, 22, 3, 0, 10) }
    @{ Row = 122; Country = 'Venezuela'; Values = @(329, 0, 142, 177, 3, 0, 10) }
    @{ Row = 134; Country = 'Islas Feroe'; Values = @(187, 0, 181, 6, 0, 0, 0) }
    @{ Row = 147; Country = 'Cabo Verde'; Values = @(109, 0, 2, 106, 0, 0, 1) }
    @{ Row = 206; Country = 'Santo Tome y Principe'; Values = @(8, 4, 4, 4, 0, 0, 0) }
    @{ Row = 207; Country = 'Butan'; Values = @(7, 0, 5, 2, 0, 0, 0) }
    @{ Row = 208; Country = 'Mauritania'; Values = @(7, 0, 6, 0, 0, 0, 1) }
    @{ Row = 209; Country = 'Sudan del Sur'; Values = @(6, 0, 0, 6, 0, 0, 0) }
    @{ Row = 210; Country = 'Islas Virgenes Britanicas'; Values = @(6, 0, 3, 2, 0, 0, 1) }
    @{ Row = 211; Country = 'Sahara Occidental'; Values = @(6, 0, 5, 1, 0, 0, 0) }
    @{ Row = 212; Country = 'San Bartolome'; Values = @(6, 0, 6, 0, 0, 0, 0) }
    @{ Row = 213; Country = 'Bonaire, San Eustaquio y Saba'; Values = @(5, 0, 0, 5, 0, 0, 0) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Country
    for ($i = 0; $i -lt $u.Values.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $u.Values[$i]
    }
}

